$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "98×95="
$t.Cell(1, 2).Range.Text = "52×53="
$t.Cell(1, 3).Range.Text = "94×88="
$t.Cell(1, 4).Range.Text = "45×96="
$t.Cell(1, 5).Range.Text = "58×17="
$t.Cell(2, 1).Range.Text = "63×61="
$t.Cell(2, 2).Range.Text = "90×76="
$t.Cell(2, 3).Range.Text = "95×34="
$t.Cell(2, 4).Range.Text = "37×13="
$t.Cell(2, 5).Range.Text = "69×83="
$t.Cell(3, 1).Range.Text = "71×95="
$t.Cell(3, 2).Range.Text = "60×43="
$t.Cell(3, 3).Range.Text = "61×41="
$t.Cell(3, 4).Range.Text = "82×49="
$t.Cell(3, 5).Range.Text = "66×83="
$t.Cell(4, 1).Range.Text = "56×29="
$t.Cell(4, 2).Range.Text = "33×60="
$t.Cell(4, 3).Range.Text = "82×80="
$t.Cell(4, 4).Range.Text = "50×31="
$t.Cell(4, 5).Range.Text = "69×51="
$t.Cell(5, 1).Range.Text = "62×40="
$t.Cell(5, 2).Range.Text = "60×60="
$t.Cell(5, 3).Range.Text = "42×69="
$t.Cell(5, 4).Range.Text = "39×81="
$t.Cell(5, 5).Range.Text = "26×52="
$t.Cell(6, 1).Range.Text = "29×48="
$t.Cell(6, 2).Range.Text = "22×28="
$t.Cell(6, 3).Range.Text = "92×99="
$t.Cell(6, 4).Range.Text = "32×33="
$t.Cell(6, 5).Range.Text = "85×49="
$t.Cell(7, 1).Range.Text = "60×26="
$t.Cell(7, 2).Range.Text = "62×41="
$t.Cell(7, 3).Range.Text = "68×28="
$t.Cell(7, 4).Range.Text = "10×39="
$t.Cell(7, 5).Range.Text = "10×11="
$t.Cell(8, 1).Range.Text = "44×90="
$t.Cell(8, 2).Range.Text = "57×95="
$t.Cell(8, 3).Range.Text = "25×14="
$t.Cell(8, 4).Range.Text = "30×36="
$t.Cell(8, 5).Range.Text = "98×50="
$t.Cell(9, 1).Range.Text = "28×34="
$t.Cell(9, 2).Range.Text = "18×62="
$t.Cell(9, 3).Range.Text = "30×87="
$t.Cell(9, 4).Range.Text = "100×39="
$t.Cell(9, 5).Range.Text = "64×53="
$t.Cell(10, 1).Range.Text = "89×97="
$t.Cell(10, 2).Range.Text = "100×22="
$t.Cell(10, 3).Range.Text = "38×13="
$t.Cell(10, 4).Range.Text = "21×26="
$t.Cell(10, 5).Range.Text = "36×75="
$t.Cell(11, 1).Range.Text = "71×74="
$t.Cell(11, 2).Range.Text = "62×90="
$t.Cell(11, 3).Range.Text = "19×39="
$t.Cell(11, 4).Range.Text = "43×95="
$t.Cell(11, 5).Range.Text = "51×62="
$t.Cell(12, 1).Range.Text = "17×11="
$t.Cell(12, 2).Range.Text = "87×27="
$t.Cell(12, 3).Range.Text = "47×42="
$t.Cell(12, 4).Range.Text = "74×30="
$t.Cell(12, 5).Range.Text = "41×74="
$t.Cell(13, 1).Range.Text = "60×70="
$t.Cell(13, 2).Range.Text = "37×83="
$t.Cell(13, 3).Range.Text = "67×37="
$t.Cell(13, 4).Range.Text = "89×10="
$t.Cell(13, 5).Range.Text = "32×76="
$t.Cell(14, 1).Range.Text = "75×100="
$t.Cell(14, 2).Range.Text = "18×65="
$t.Cell(14, 3).Range.Text = "34×78="
$t.Cell(14, 4).Range.Text = "51×41="
$t.Cell(14, 5).Range.Text = "54×30="
$t.Cell(15, 1).Range.Text = "13×95="
$t.Cell(15, 2).Range.Text = "18×40="
$t.Cell(15, 3).Range.Text = "30×46="
$t.Cell(15, 4).Range.Text = "29×76="
$t.Cell(15, 5).Range.Text = "87×49="
$t.Cell(16, 1).Range.Text = "59×30="
$t.Cell(16, 2).Range.Text = "39×45="
$t.Cell(16, 3).Range.Text = "12×73="
$t.Cell(16, 4).Range.Text = "88×48="
$t.Cell(16, 5).Range.Text = "80×58="
$t.Cell(17, 1).Range.Text = "72×42="
$t.Cell(17, 2).Range.Text = "91×25="
$t.Cell(17, 3).Range.Text = "28×36="
$t.Cell(17, 4).Range.Text = "84×20="
$t.Cell(17, 5).Range.Text = "82×64="
$t.Cell(18, 1).Range.Text = "15×22="
$t.Cell(18, 2).Range.Text = "68×37="
$t.Cell(18, 3).Range.Text = "65×36="
$t.Cell(18, 4).Range.Text = "50×27="
$t.Cell(18, 5).Range.Text = "55×59="
$t.Cell(19, 1).Range.Text = "71×46="
$t.Cell(19, 2).Range.Text = "98×69="
$t.Cell(19, 3).Range.Text = "11×15="
$t.Cell(19, 4).Range.Text = "63×23="
$t.Cell(19, 5).Range.Text = "84×26="
$t.Cell(20, 1).Range.Text = "35×97="
$t.Cell(20, 2).Range.Text = "25×39="
$t.Cell(20, 3).Range.Text = "58×11="
$t.Cell(20, 4).Range.Text = "37×16="
$t.Cell(20, 5).Range.Text = "94×37="
